{"js": "// Fix typo in the variables worksheet:\n// \"Write a program the says the largest of three numbers given by a user\"\n// should read\n// \"Write a program that says the largest of three numbers given by a user\"\nconst body = context.document.body;\n\n// The mis-typed phrase \"the \" (should be \"that \") is unique in the document -\n// every other instance already reads \"Write a program that ...\".\nconst results = body.search(\"Write a program the \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(`Expected exactly 1 match for the typo, found ${results.items.length}`);\n}\n\n// Replace just the matched range's text in place, preserving its formatting\n// (font, color, etc.) and leaving every other run in the paragraph untouched.\nresults.items[0].insertText(\"Write a program that \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix typo in the variables worksheet:\n# \"Write a program the says the largest of three numbers given by a user\"\n# should read\n# \"Write a program that says the largest of three numbers given by a user\"\n$d = $word.ActiveDocument\n\n# The mis-typed phrase \"the \" (should be \"that \") only occurs once in the\n# document - every other \"Write a program ...\" sentence already says \"that \".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"Write a program the \",  # FindText\n    $true,                    # MatchCase\n    $false,                   # MatchWholeWord\n    $false,                   # MatchWildcards\n    $false,                   # MatchSoundsLike\n    $false,                   # MatchAllWordForms\n    $true,                    # Forward\n    1,                        # Wrap (wdFindContinue)\n    $false,                   # Format\n    \"Write a program that \",  # ReplaceWith\n    1                         # Replace (wdReplaceOne)\n)\n\nif (-not $found) {\n    throw \"Typo text 'Write a program the ' was not found\"\n}\n"}
